$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: first "Senior Software Engineer" heading (Cotiviti entry,
# paragraph 17) -> "Advisor Product Developer"
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
if ($p17.Range.Text -ne "Senior Software Engineer`r") {
    Write-Host "WARN: paragraph 17 text unexpected: [$($p17.Range.Text)]"
}
$r17 = $d.Range($p17.Range.Start, $p17.Range.End - 1)
$r17.Text = "Advisor Product Developer"

# ---------------------------------------------------------------------
# Change 2: "retrieval management" -> "Managed Care" (paragraph 18,
# the bullet right after the heading changed above)
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs(18)
$scope18 = $d.Range($p18.Range.Start, $p18.Range.End)
$found18 = $scope18.Find.Execute("retrieval management", $true, $false, $false, $false, $false, $true, 1, $false, "Managed Care", 2)
if (-not $found18) {
    Write-Host "WARN: 'retrieval management' not found in paragraph 18"
}

# ---------------------------------------------------------------------
# Change 3: "MS SQL server, Oracle" -> "MS SQL Server, Oracle"
# (paragraph 19, capitalising the 's' in "server")
# ---------------------------------------------------------------------
$p19 = $d.Paragraphs(19)
$scope19 = $d.Range($p19.Range.Start, $p19.Range.End)
$found19 = $scope19.Find.Execute("MS SQL server, Oracle", $true, $false, $false, $false, $false, $true, 1, $false, "MS SQL Server, Oracle", 2)
if (-not $found19) {
    Write-Host "WARN: 'MS SQL server, Oracle' not found in paragraph 19"
}

# ---------------------------------------------------------------------
# Change 4: "ASP.NET MVC, MS SQL server" -> "ASP.NET REST, MS SQL Server"
# (paragraph 31). Word leaves a "_GoBack" bookmark at the spot of the
# most recent edit (the "server" -> "Server" capitalisation), sitting
# between the "S" and "erver". We reproduce this exactly:
#   1. locate "server" and drop the bookmark between its 1st and 2nd
#      characters (this also keeps the trailing ", Node.js" run intact
#      across the later edits, since a bookmark can't be merged across),
#   2. do the "ASP.NET MVC" -> "ASP.NET REST" replacement,
#   3. capitalise the "s" immediately before the bookmark to "S".
# ---------------------------------------------------------------------
$p31 = $d.Paragraphs(31)
$scope31 = $d.Range($p31.Range.Start, $p31.Range.End)
$foundServer = $scope31.Find.Execute("server")
if ($foundServer) {
    $sPos = $scope31.Start
    $bmPoint = $d.Range($sPos + 1, $sPos + 1)
    $d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

    $p31b = $d.Paragraphs(31)
    $scope31b = $d.Range($p31b.Range.Start, $p31b.Range.End)
    $foundMvc = $scope31b.Find.Execute("ASP.NET MVC", $true, $false, $false, $false, $false, $true, 1, $false, "ASP.NET REST", 2)
    if (-not $foundMvc) {
        Write-Host "WARN: 'ASP.NET MVC' not found in paragraph 31"
    }

    $bm = $d.Bookmarks("_GoBack")
    $bmRange = $bm.Range
    $sCharRange = $d.Range($bmRange.Start - 1, $bmRange.Start)
    if ($sCharRange.Text -ne "s") {
        Write-Host "WARN: char before bookmark is [$($sCharRange.Text)], expected 's'"
    }
    $sCharRange.Text = "S"
} else {
    Write-Host "WARN: 'server' not found in paragraph 31"
}

Write-Host "Para17: $($d.Paragraphs(17).Range.Text)"
Write-Host "Para18: $($d.Paragraphs(18).Range.Text)"
Write-Host "Para19: $($d.Paragraphs(19).Range.Text)"
Write-Host "Para31: $($d.Paragraphs(31).Range.Text)"
